# Applies the "cryptos list" price/volume update described by the commit
# "Updated cryptos list on Thu Apr 11 13:50:40 UTC 2024 with GitHub Actions".
#
# For most rows only the Price (D) and Volume(1h) (E) columns move. A few
# coins swapped rank (and therefore swapped rows) between snapshots, so for
# those row-pairs the Coin (B), Link (C), Price (D) and Volume(1h) (E)
# columns are all rewritten to reflect the new row contents.
#
# All touched cells originally hold plain text (prices such as "0.580" or
# "3.354.82" are not real numbers - they are pre-formatted display strings).
# Force the Text number format before each write so Excel's COM layer does
# not "helpfully" reinterpret them as numbers and drop trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.866.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  +3.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.504.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  +1.89%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "600.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  +3.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "172.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  +4.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.608"
$ws.Range("E7").NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  +1.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.498.21"
$ws.Range("E8").NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  +1.96%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.192"
$ws.Range("E10").NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  +4.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.26"
$ws.Range("E11").NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  +8.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.580"
$ws.Range("E12").NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  +2.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "45.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  +0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000273"
$ws.Range("E14").NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  +1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "4.084.13"
$ws.Range("E15").NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  +2.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "8.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "  +0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "602.89"
$ws.Range("E17").NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.536.13"
$ws.Range("E18").NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  +2.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "69.998.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  +3.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.119"
$ws.Range("E20").NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  +1.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  +0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.867"
$ws.Range("E22").NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  +0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  -15.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "15.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  +1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "95.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  +0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.70"
$ws.Range("E26").NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  -0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  +0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "33.65"
$ws.Range("E29").NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  +5.37%  "

# Row 30/31 swapped: RenderToken now ranks above Bittensor
$ws.Range("B30").NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.95"
$ws.Range("E30").NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  +0.23%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "Bittensor"
$ws.Range("C31").NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "720.99"
$ws.Range("E31").NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  +23.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "8.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  -2.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  +3.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.28"
$ws.Range("E35").NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  +0.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0995"
$ws.Range("E36").NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  -0.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.54"
$ws.Range("E37").NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  +5.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "10.66"
$ws.Range("E38").NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  +0.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0473"
$ws.Range("E39").NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  +10.38%  "

# Row 40/41 swapped: FirstDigitalUSD now ranks above OKB
$ws.Range("B40").NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.00"

$ws.Range("B41").NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "OKB"
$ws.Range("C41").NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "56.46"
$ws.Range("E41").NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.141"
$ws.Range("E42").NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  +5.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.340.86"
$ws.Range("E43").NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  -0.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.313"
$ws.Range("E44").NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  -1.25%  "

# Row 45/46 swapped: InjectiveProtocol now ranks above PEPE
$ws.Range("B45").NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = "InjectiveProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "32.24"
$ws.Range("E45").NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  +0.17%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = "PEPE"
$ws.Range("C46").NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0₃0689"
$ws.Range("E46").NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  +2.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.89"
$ws.Range("E47").NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  +6.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  +3.60%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "  +1.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "133.30"
$ws.Range("E50").NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  +1.12%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  -0.04%  "
